# StructureDefinition-FrRatioUcum.xlsx — refresh generated IG export
# (ci-build message 40e73818b405ab3cd55ea5bdf7793eae299a4c02)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Metadata" sheet: canonical URL moved, publication date refreshed,
# and the stale copyright note dropped.
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "https://hl7.fr/fhir/fr/medication/StructureDefinition/FrRatioUcum"
$meta.Range("B8").Value = "2025-04-10T15:35:36+00:00"
$meta.Range("B14").Value = ""

# ---------------------------------------------------------------------
# "Elements" sheet: Ratio.numerator (row 5) and Ratio.denominator
# (row 6) get marked as summary elements, their verbose comments /
# constraints / v2 mappings are trimmed, and the RIM mapping is
# replaced with the short ".numerator" / ".denominator" path.
# ---------------------------------------------------------------------
$el = $wb.Worksheets.Item("Elements")

$numeratorRow = 5
$denominatorRow = 6

foreach ($row in @($numeratorRow, $denominatorRow)) {
    $el.Cells.Item($row, 10).Value = "Y"     # J: Is Summary?
    $el.Cells.Item($row, 14).Value = ""      # N: Comments
    $el.Cells.Item($row, 35).Value = ""      # AI: Condition(s)
    $el.Cells.Item($row, 36).Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`n"  # AJ: Constraint(s)
    $el.Cells.Item($row, 37).Value = ""      # AK: Mapping: HL7 v2 Mapping
}

$el.Range("AL5").Value = ".numerator"        # AL: Mapping: RIM Mapping
$el.Range("AL6").Value = ".denominator"

# Column AL narrows to fit the now much shorter RIM mapping text.
$el.Columns.Item(38).ColumnWidth = 24.2
